$d = $word.ActiveDocument

# --- Step 1: capture (via the clipboard) the "Meta description" paragraph,
# which already has the run layout/formatting (leading empty run + a bold
# run) that the new end-of-document heading line should reuse. ---
$metaPara = $d.Paragraphs.Item(2)
$metaTemplateRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaTemplateRange.Copy()

# --- Step 2: remove the "Meta description" paragraph that currently follows
# the title heading. ---
$metaPara.Range.Delete()

# --- Step 3: paste the captured paragraph right before the final
# (image-prompt) paragraph, then turn its text into the bold
# "Play Dollar Bomb Free: Slot Game Review" heading line. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.Paste()

$newPara = $d.Paragraphs.Item($count)
$newRange = $newPara.Range
$newRange.Find.Execute("Meta description: Discover more about Dollar Bomb by CQ9 Gaming in our review. Play Dollar Bomb for free online and experience great winning potential with an eastern theme.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Play Dollar Bomb Free: Slot Game Review", 2)

# --- Step 4: update the final paragraph's text (keeping its existing italic
# run formatting) from the old image-generation prompt to the new copy. ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.Find.Execute("Create a feature image for Dollar Bomb with the following specifications: Design a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be holding a bag of golden coins and standing in front of a grand castle. The castle should be in the background, with Dollar Bomb's name prominently displayed above it in bold, golden letters. Use bright, eye-catching colors to make the image stand out and add a touch of whimsy to the design. The image should convey excitement, joy, and the promise of big wins to players.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Discover more about Dollar Bomb by CQ9 Gaming in our review. Play Dollar Bomb for free online and experience great winning potential with an eastern theme.", 2)
